$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3788.1267
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 3868.9421
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 11606.8263
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -11942.8263
$ws.Range("H58").Value = 3628.4707
$ws.Range("I58").Value = 108.4
$ws.Range("J58").Value = 8657.143
$ws.Range("K58").Value = 325.2
$ws.Range("L58").Value = 25971.429
$ws.Range("M58").Value = -175.2
$ws.Range("N58").Value = -26271.429
$ws.Range("H64").Value = 4036089.8
$ws.Range("I64").Value = 8931991
$ws.Range("J64").Value = 4170.7646
$ws.Range("K64").Value = 8931991
$ws.Range("L64").Value = 4170.7646
$ws.Range("M64").Value = -8931743
$ws.Range("N64").Value = -4666.7646
$ws.Range("H67").Value = 4036089.8
$ws.Range("I67").Value = 8931991
$ws.Range("J67").Value = 4170.7646
$ws.Range("K67").Value = 8931991
$ws.Range("L67").Value = 4170.7646
$ws.Range("M67").Value = -8931133
$ws.Range("N67").Value = -5886.7646
$ws.Range("H98").Value = 431975.5
$ws.Range("I98").Value = 590271.1
$ws.Range("J98").Value = 2316
$ws.Range("K98").Value = 590271.1
$ws.Range("L98").Value = 2316
$ws.Range("M98").Value = -588773.1
$ws.Range("N98").Value = -5312
$ws.Range("H99").Value = 9091389
$ws.Range("I99").Value = 11363987
$ws.Range("J99").Value = 1000
$ws.Range("K99").Value = 34091961
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -34090463
$ws.Range("N99").Value = -5996
$ws.Range("H118").Value = 800
$ws.Range("I118").Value = 800
$ws.Range("K118").Value = 2400
$ws.Range("M118").Value = -743
$ws.Range("H121").Value = 612.2727
$ws.Range("J121").Value = 612.2727
$ws.Range("L121").Value = 1836.8181
$ws.Range("N121").Value = -5330.8181
$ws.Range("H122").Value = 431975.5
$ws.Range("I122").Value = 590271.1
$ws.Range("J122").Value = 2316
$ws.Range("K122").Value = 1770813.3
$ws.Range("L122").Value = 6948
$ws.Range("M122").Value = -1768363.3
$ws.Range("N122").Value = -11848
$ws.Range("H135").Value = 3842.2
$ws.Range("I135").Value = 3065
$ws.Range("J135").Value = 5008
$ws.Range("K135").Value = 27585
$ws.Range("L135").Value = 45072
$ws.Range("M135").Value = -25050
$ws.Range("N135").Value = -50142
$ws.Range("H138").Value = 2652.23
$ws.Range("I138").Value = 1553.12
$ws.Range("J138").Value = 3018.6
$ws.Range("K138").Value = 4659.36
$ws.Range("L138").Value = 9055.799999999999
$ws.Range("M138").Value = 480.6400000000003
$ws.Range("N138").Value = -19335.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 41711.977
$ws.Range("I32").Value = 8801.546
$ws.Range("J32").Value = 196861.14
$ws.Range("K32").Value = 8801.546
$ws.Range("L32").Value = 196861.14
$ws.Range("M32").Value = -8514.546
$ws.Range("N32").Value = -197435.14
$ws.Range("H63").Value = 14320
$ws.Range("I63").Value = 17333.334
$ws.Range("K63").Value = 17333.334
$ws.Range("M63").Value = -16647.334
$ws.Range("H66").Value = 14320
$ws.Range("I66").Value = 17333.334
$ws.Range("K66").Value = 86666.67
$ws.Range("M66").Value = -83234.67
$ws.Range("H74").Value = 4192.6665
$ws.Range("I74").Value = 923.7037
$ws.Range("J74").Value = 10076.8
$ws.Range("K74").Value = 923.7037
$ws.Range("L74").Value = 10076.8
$ws.Range("M74").Value = -49.70370000000003
$ws.Range("N74").Value = -11824.8
$ws.Range("H77").Value = 4192.6665
$ws.Range("I77").Value = 923.7037
$ws.Range("J77").Value = 10076.8
$ws.Range("K77").Value = 4618.5185
$ws.Range("L77").Value = 50384
$ws.Range("M77").Value = -250.5185000000001
$ws.Range("N77").Value = -59120
$ws.Range("H132").Value = 2664.923
$ws.Range("I132").Value = 2015.3334
$ws.Range("J132").Value = 5393.2
$ws.Range("K132").Value = 6046.0002
$ws.Range("L132").Value = 16179.6
$ws.Range("M132").Value = -3516.0002
$ws.Range("N132").Value = -21239.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2426.9092
$ws.Range("I20").Value = 2274.5
$ws.Range("J20").Value = 2833.3333
$ws.Range("K20").Value = 2274.5
$ws.Range("L20").Value = 2833.3333
$ws.Range("M20").Value = -2027.5
$ws.Range("N20").Value = -3327.3333
$ws.Range("H86").Value = 1662.3077
$ws.Range("I86").Value = 1799.2
$ws.Range("J86").Value = 1576.75
$ws.Range("K86").Value = 1799.2
$ws.Range("L86").Value = 1576.75
$ws.Range("M86").Value = -676.2
$ws.Range("N86").Value = -3822.75
$ws.Range("H89").Value = 1662.3077
$ws.Range("I89").Value = 1799.2
$ws.Range("J89").Value = 1576.75
$ws.Range("K89").Value = 8996
$ws.Range("L89").Value = 7883.75
$ws.Range("M89").Value = -3380
$ws.Range("N89").Value = -19115.75
$ws.Range("H134").Value = 2254.578
$ws.Range("I134").Value = 1829.4736
$ws.Range("J134").Value = 4562.2856
$ws.Range("K134").Value = 5488.4208
$ws.Range("L134").Value = 13686.8568
$ws.Range("M134").Value = -2953.4208
$ws.Range("N134").Value = -18756.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3984.3647
$ws.Range("I31").Value = 2233.3635
$ws.Range("J31").Value = 6552.5
$ws.Range("K31").Value = 2233.3635
$ws.Range("L31").Value = 6552.5
$ws.Range("M31").Value = -1938.3635
$ws.Range("N31").Value = -7142.5
$ws.Range("H34").Value = 3984.3647
$ws.Range("I34").Value = 2233.3635
$ws.Range("J34").Value = 6552.5
$ws.Range("K34").Value = 2233.3635
$ws.Range("L34").Value = 6552.5
$ws.Range("M34").Value = -2031.3635
$ws.Range("N34").Value = -6956.5
$ws.Range("H58").Value = 20001776
$ws.Range("I58").Value = 27028160
$ws.Range("J58").Value = 3609.6155
$ws.Range("K58").Value = 27028160
$ws.Range("L58").Value = 3609.6155
$ws.Range("M58").Value = -27027957
$ws.Range("N58").Value = -4015.6155
$ws.Range("H99").Value = 13892642
$ws.Range("I99").Value = 3933.0588
$ws.Range("K99").Value = 3933.0588
$ws.Range("M99").Value = -2435.0588
$ws.Range("H126").Value = 13892642
$ws.Range("I126").Value = 3933.0588
$ws.Range("K126").Value = 11799.1764
$ws.Range("M126").Value = -9329.1764
$ws.Range("H134").Value = 31917332
$ws.Range("I134").Value = 41667924
$ws.Range("J134").Value = 21742800
$ws.Range("K134").Value = 125003772
$ws.Range("L134").Value = 65228400
$ws.Range("M134").Value = -125001237
$ws.Range("N134").Value = -65233470
$ws.Range("H136").Value = 20001776
$ws.Range("I136").Value = 27028160
$ws.Range("J136").Value = 3609.6155
$ws.Range("K136").Value = 81084480
$ws.Range("L136").Value = 10828.8465
$ws.Range("M136").Value = -81081930
$ws.Range("N136").Value = -15928.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 9181.954
$ws.Range("J39").Value = 9571.619000000001
$ws.Range("L39").Value = 28714.857
$ws.Range("N39").Value = -29302.857
$ws.Range("H55").Value = 3636.3635
$ws.Range("J55").Value = 4750
$ws.Range("L55").Value = 14250
$ws.Range("N55").Value = -14604
$ws.Range("H107").Value = 638.3077
$ws.Range("J107").Value = 700
$ws.Range("L107").Value = 2100
$ws.Range("N107").Value = -5940
$ws.Range("H131").Value = 14494938
$ws.Range("I131").Value = 910
$ws.Range("K131").Value = 2730
$ws.Range("M131").Value = 2310
$ws.Range("H139").Value = 9909.091
$ws.Range("I139").Value = 9000
$ws.Range("J139").Value = 10000
$ws.Range("K139").Value = 27000
$ws.Range("L139").Value = 30000
$ws.Range("M139").Value = -21860
$ws.Range("N139").Value = -40280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6102.0625
$ws.Range("I102").Value = 4013.3
$ws.Range("J102").Value = 9583.333000000001
$ws.Range("K102").Value = 4013.3
$ws.Range("L102").Value = 9583.333000000001
$ws.Range("M102").Value = -2391.3
$ws.Range("N102").Value = -12827.333
$ws.Range("H113").Value = 1350
$ws.Range("I113").Value = 1350
$ws.Range("K113").Value = 1350
$ws.Range("M113").Value = 820
$ws.Range("H132").Value = 3264.9644
$ws.Range("I132").Value = 3036.4
$ws.Range("J132").Value = 3836.375
$ws.Range("K132").Value = 9109.200000000001
$ws.Range("L132").Value = 11509.125
$ws.Range("M132").Value = -6579.200000000001
$ws.Range("N132").Value = -16569.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3384.7144
$ws.Range("I122").Value = 2566.889
$ws.Range("J122").Value = 3772.1052
$ws.Range("K122").Value = 7700.667
$ws.Range("L122").Value = 11316.3156
$ws.Range("M122").Value = -5250.667
$ws.Range("N122").Value = -16216.3156
$ws.Range("H132").Value = 5181.55
$ws.Range("I132").Value = 3480.2727
$ws.Range("J132").Value = 7260.8887
$ws.Range("K132").Value = 10440.8181
$ws.Range("L132").Value = 21782.6661
$ws.Range("M132").Value = -7910.8181
$ws.Range("N132").Value = -26842.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3465.4783
$ws.Range("I132").Value = 3506.7144
$ws.Range("J132").Value = 3401.3333
$ws.Range("K132").Value = 10520.1432
$ws.Range("L132").Value = 10203.9999
$ws.Range("M132").Value = -7990.143199999999
$ws.Range("N132").Value = -15263.9999
$ws.Range("H136").Value = 2659.5117
$ws.Range("I136").Value = 719.129
$ws.Range("J136").Value = 7672.1665
$ws.Range("K136").Value = 2157.387
$ws.Range("L136").Value = 23016.4995
$ws.Range("M136").Value = 392.6129999999998
$ws.Range("N136").Value = -28116.4995
